$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.044118602376014
$ws.Range("D2").Value = 1.052413044249146
$ws.Range("E2").Value = 1.042015255067438
$ws.Range("F2").Value = 1.06129037144762
$ws.Range("I2").Value = 1.045797667916304
$ws.Range("J2").Value = 1.049185399982641
$ws.Range("K2").Value = 1.055161846104649
$ws.Range("L2").Value = 1.044793137517146
$ws.Range("M2").Value = 1.064014833492214
$ws.Range("N2").Value = 1.020185035266913

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.045195091153771
$ws.Range("D3").Value = 1.053318796776867
$ws.Range("E3").Value = 1.04293628443152
$ws.Range("F3").Value = 1.062429079701077
$ws.Range("I3").Value = 1.046159415554812
$ws.Range("J3").Value = 1.049908384772614
$ws.Range("K3").Value = 1.055880200810739
$ws.Range("L3").Value = 1.045524623071276
$ws.Range("M3").Value = 1.064967316977219
$ws.Range("N3").Value = 1.020433315968315

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.045891626449108
$ws.Range("D4").Value = 1.053904836583693
$ws.Range("E4").Value = 1.043532525358555
$ws.Range("F4").Value = 1.063166344919341
$ws.Range("I4").Value = 1.046392158420051
$ws.Range("J4").Value = 1.050375605218781
$ws.Range("K4").Value = 1.056344349391665
$ws.Range("L4").Value = 1.04599758914073
$ws.Range("M4").Value = 1.065583495807906
$ws.Range("N4").Value = 1.020593581147256

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.046184444423097
$ws.Range("D5").Value = 1.05415119739637
$ws.Range("E5").Value = 1.043783250302278
$ws.Range("F5").Value = 1.063476397929887
$ws.Range("I5").Value = 1.046489684656141
$ws.Range("J5").Value = 1.050571881341492
$ws.Range("K5").Value = 1.056539315896045
$ws.Range("L5").Value = 1.046196339285977
$ws.Range("M5").Value = 1.065842503618876
$ws.Range("N5").Value = 1.020660863404083

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.04623360954369
$ws.Range("D6").Value = 1.054192561852303
$ws.Range("E6").Value = 1.043825351945115
$ws.Range("F6").Value = 1.063528463458763
$ws.Range("I6").Value = 1.046506041042861
$ws.Range("J6").Value = 1.050604828561732
$ws.Range("K6").Value = 1.056572042156886
$ws.Range("L6").Value = 1.04622970533655
$ws.Range("M6").Value = 1.06588599019412
$ws.Range("N6").Value = 1.020672154921788

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.045895539118036
$ws.Range("D7").Value = 1.053908128509627
$ws.Range("E7").Value = 1.043535875301327
$ws.Range("F7").Value = 1.063170487443252
$ws.Range("I7").Value = 1.046393462823269
$ws.Range("J7").Value = 1.050378228433135
$ws.Range("K7").Value = 1.056346955178421
$ws.Range("L7").Value = 1.046000245183117
$ws.Range("M7").Value = 1.065586956815452
$ws.Range("N7").Value = 1.02059448054248

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.044482412262068
$ws.Range("D8").Value = 1.052719156318909
$ws.Range("E8").Value = 1.042326464548004
$ws.Range("F8").Value = 1.061675111258432
$ws.Range("I8").Value = 1.045920198051401
$ws.Range("J8").Value = 1.049429860154028
$ws.Range("K8").Value = 1.055404757208758
$ws.Range("L8").Value = 1.045040419756818
$ws.Range("M8").Value = 1.064336759169036
$ws.Range("N8").Value = 1.020269023392067

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.041992088415384
$ws.Range("D9").Value = 1.050623713274586
$ws.Range("E9").Value = 1.040197429248368
$ws.Range("F9").Value = 1.059043464021913
$ws.Range("I9").Value = 1.045076041565082
$ws.Range("J9").Value = 1.047754122787067
$ws.Range("K9").Value = 1.053739313134801
$ws.Range("L9").Value = 1.043346374995791
$ws.Range("M9").Value = 1.062132659050097
$ws.Range("N9").Value = 1.019692550505427

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.040331694072561
$ws.Range("D10").Value = 1.049226537881583
$ws.Range("E10").Value = 1.03877949244934
$ws.Range("F10").Value = 1.057291298833279
$ws.Range("I10").Value = 1.044506406486332
$ws.Range("J10").Value = 1.046633870680484
$ws.Range("K10").Value = 1.052625535790538
$ws.Range("L10").Value = 1.042215189530505
$ws.Range("M10").Value = 1.060662511946408
$ws.Range("N10").Value = 1.019306238709303

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.039612674072357
$ws.Range("D11").Value = 1.048621494667264
$ws.Range("E11").Value = 1.038165847709126
$ws.Range("F11").Value = 1.056533124648735
$ws.Range("I11").Value = 1.04425811876756
$ws.Range("J11").Value = 1.046148053011116
$ws.Range("K11").Value = 1.052142430347424
$ws.Range("L11").Value = 1.041724940633298
$ws.Range("M11").Value = 1.060025741100852
$ws.Range("N11").Value = 1.019138488497623

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.03934558842447
$ws.Range("D12").Value = 1.048396745977135
$ws.Range("E12").Value = 1.037937962455102
$ws.Range("F12").Value = 1.056251583330807
$ws.Range("I12").Value = 1.044165648275389
$ws.Range("J12").Value = 1.04596748705721
$ws.Range("K12").Value = 1.051962858143196
$ws.Range("L12").Value = 1.041542774370787
$ws.Range("M12").Value = 1.059789187464933
$ws.Range("N12").Value = 1.019076107299835

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.039402879650842
$ws.Range("D13").Value = 1.048444955729916
$ws.Range("E13").Value = 1.037986842366178
$ws.Range("F13").Value = 1.056311971363018
$ws.Range("I13").Value = 1.044185494616838
$ws.Range("J13").Value = 1.046006224123987
$ws.Range("K13").Value = 1.052001382666271
$ws.Range("L13").Value = 1.041581852641561
$ws.Range("M13").Value = 1.059839930323561
$ws.Range("N13").Value = 1.019089491504813

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.039590596890226
$ws.Range("D14").Value = 1.048602917048673
$ws.Range("E14").Value = 1.038147009628806
$ws.Range("F14").Value = 1.056509850760676
$ws.Range("I14").Value = 1.044250480134589
$ws.Range("J14").Value = 1.046133129652182
$ws.Range("K14").Value = 1.052127589411239
$ws.Range("L14").Value = 1.041709884060448
$ws.Range("M14").Value = 1.060006188089256
$ws.Range("N14").Value = 1.019133333502788

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.039706254378737
$ws.Range("D15").Value = 1.048700241073715
$ws.Range("E15").Value = 1.038245700553297
$ws.Range("F15").Value = 1.056631781152354
$ws.Range("I15").Value = 1.044290487336093
$ws.Range("J15").Value = 1.046211305534217
$ws.Range("K15").Value = 1.052205332936692
$ws.Range("L15").Value = 1.041788759693298
$ws.Range("M15").Value = 1.060108621184396
$ws.Range("N15").Value = 1.0191603365571

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.040379411761375
$ws.Range("D16").Value = 1.049266691394028
$ws.Range("E16").Value = 1.038820225005283
$ws.Range("F16").Value = 1.057341627370871
$ws.Range("I16").Value = 1.044522850124236
$ws.Range("J16").Value = 1.04666609715945
$ws.Range("K16").Value = 1.052657580345792
$ws.Range("L16").Value = 1.042247716489861
$ws.Range("M16").Value = 1.060704768344255
$ws.Range("N16").Value = 1.019317361738921

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.040801649397568
$ws.Range("D17").Value = 1.04962199534964
$ws.Range("E17").Value = 1.03918069767144
$ws.Range("F17").Value = 1.0577870353082
$ws.Range("I17").Value = 1.04466816810069
$ws.Range("J17").Value = 1.04695117730387
$ws.Range("K17").Value = 1.052941040130357
$ws.Range("L17").Value = 1.042535490543967
$ws.Range("M17").Value = 1.061078665219072
$ws.Range("N17").Value = 1.01941573241645

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.041047927899418
$ws.Range("D18").Value = 1.049829232640028
$ws.Range("E18").Value = 1.039390987242633
$ws.Range("F18").Value = 1.058046884903487
$ws.Range("I18").Value = 1.044752772169161
$ws.Range("J18").Value = 1.047117388241977
$ws.Range("K18").Value = 1.053106297119753
$ws.Range("L18").Value = 1.042703302019584
$ws.Range("M18").Value = 1.06129673488308
$ws.Range("N18").Value = 1.019473064570292

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.041131901585578
$ws.Range("D19").Value = 1.049899894272331
$ws.Range("E19").Value = 1.039462695974896
$ws.Range("F19").Value = 1.058135495512866
$ws.Range("I19").Value = 1.044781593270026
$ws.Range("J19").Value = 1.047174049760187
$ws.Range("K19").Value = 1.053162631868878
$ws.Range("L19").Value = 1.042760514233046
$ws.Range("M19").Value = 1.061371087953338
$ws.Range("N19").Value = 1.019492605584663

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.040756347877063
$ws.Range("D20").Value = 1.049583875167023
$ws.Range("E20").Value = 1.03914201905987
$ws.Range("F20").Value = 1.0577392420133
$ws.Range("I20").Value = 1.044652593144063
$ws.Range("J20").Value = 1.046920598298604
$ws.Range("K20").Value = 1.05291063588641
$ws.Range("L20").Value = 1.042504619483847
$ws.Range("M20").Value = 1.061038551493637
$ws.Range("N20").Value = 1.019405182910712

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.03953531911391
$ws.Range("D21").Value = 1.048556401618611
$ws.Range("E21").Value = 1.038099842989708
$ws.Range("F21").Value = 1.056451578051979
$ws.Range("I21").Value = 1.044231350294523
$ws.Range("J21").Value = 1.046095762221366
$ws.Range("K21").Value = 1.052090428132324
$ws.Range("L21").Value = 1.041672183821801
$ws.Range("M21").Value = 1.059957230113841
$ws.Range("N21").Value = 1.019120425095256

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.038767553875976
$ws.Range("D22").Value = 1.047910337813799
$ws.Range("E22").Value = 1.037444873203363
$ws.Range("F22").Value = 1.055642425555574
$ws.Range("I22").Value = 1.043965078417523
$ws.Range("J22").Value = 1.045576508978478
$ws.Range("K22").Value = 1.051574005701461
$ws.Range("L22").Value = 1.041148416578326
$ws.Range("M22").Value = 1.059277194810416
$ws.Range("N22").Value = 1.018940974040137

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.039174566345141
$ws.Range("D23").Value = 1.048252833281973
$ws.Range("E23").Value = 1.037792057825546
$ws.Range("F23").Value = 1.056071329846243
$ws.Range("I23").Value = 1.044106368820348
$ws.Range("J23").Value = 1.045851836285319
$ws.Range("K23").Value = 1.051847839835688
$ws.Range("L23").Value = 1.041426111760938
$ws.Range("M23").Value = 1.059637710234974
$ws.Range("N23").Value = 1.019036143537549

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.040776817701102
$ws.Range("D24").Value = 1.049601100053511
$ws.Range("E24").Value = 1.039159496160643
$ws.Range("F24").Value = 1.057760837589187
$ws.Range("I24").Value = 1.04465963128292
$ws.Range("J24").Value = 1.046934415855539
$ws.Range("K24").Value = 1.052924374503485
$ws.Range("L24").Value = 1.042518568917523
$ws.Range("M24").Value = 1.061056677215064
$ws.Range("N24").Value = 1.019409949919725

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.042635925737189
$ws.Range("D25").Value = 1.051165473366923
$ws.Range("E25").Value = 1.040747586154485
$ws.Range("F25").Value = 1.05972340718711
$ws.Range("I25").Value = 1.045295485815327
$ws.Range("J25").Value = 1.048187885869598
$ws.Range("K25").Value = 1.054170483605589
$ws.Range("L25").Value = 1.043784647331742
$ws.Range("M25").Value = 1.062702602498928
$ws.Range("N25").Value = 1.019841934777432
